$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.267.01"
$ws.Range("E2").Value = "  -2.79%  "

$ws.Range("D3").Value = "2.208.34"
$ws.Range("E3").Value = "  -2.95%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "107.23"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -12.90%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "296.46"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +11.36%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.623"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -2.66%  "

$ws.Range("E8").Value = "  -0.23%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.593"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -5.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.44"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -9.77%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0907"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -4.18%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.28"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.71"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -6.88%  "

$ws.Range("E14").Value = "  -3.24%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.941"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +4.39%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.87"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.51%  "

$ws.Range("D17").Value = "2.538.47"
$ws.Range("E17").Value = "  -3.23%  "

$ws.Range("D18").Value = "2.231.85"
$ws.Range("E18").Value = "  -1.88%  "

$ws.Range("D19").Value = "41.978.23"
$ws.Range("E19").Value = "  -3.81%  "

$ws.Range("E20").Value = "  +4.56%  "

$ws.Range("E21").Value = "  -5.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.15"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.16%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.48"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +20.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.26"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -6.62%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "226.90"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.41%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.93"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -5.70%  "

$ws.Range("E27").Value = "  -1.48%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.52"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.54%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.96"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.42%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.24"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.06%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "38.01"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -9.86%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.20"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.99%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "172.54"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.79"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.85%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0874"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.64%  "

$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.00"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +8.40%  "

$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.50"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.16%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.26"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.07%  "

$ws.Range("E39").Value = "  -3.61%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0359"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.66%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.101"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.24%  "

$ws.Range("E42").Value = "  -4.43%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.42"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.47%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.228"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -4.11%  "

$ws.Range("E45").Value = "  -0.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.52"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -9.98%  "

$ws.Range("E47").Value = "  -6.00%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.38"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.47%  "

$ws.Range("E49").Value = "  +3.94%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "102.00"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.27%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.42"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.67%  "
